# Refresh the cryptocurrency price ("Price", column D) and 1-hour volume
# change ("Volume(1h)", column E) figures on Sheet1, row by row, to match
# the Jan 22 2024 data pull.
#
# Column D values that look like plain decimals (e.g. "312.01") are written
# with a leading apostrophe so Excel keeps them as text -- matching how the
# source data stores these cells (plain/European-style numeric strings,
# some of which, like "40.639.17", are not valid numbers at all) -- instead
# of silently reinterpreting them as numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '40.639.17'
$ws.Range("E2").Value = '  -2.73%  '
$ws.Range("D3").Value = '2.377.82'
$ws.Range("E3").Value = '  -4.11%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'312.01"
$ws.Range("E5").Value = '  -2.24%  '
$ws.Range("D6").Value = "'87.32"
$ws.Range("E6").Value = '  -6.79%  '
$ws.Range("D7").Value = "'0.529"
$ws.Range("E7").Value = '  -4.63%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  -5.64%  '
$ws.Range("D10").Value = "'0.0822"
$ws.Range("E10").Value = '  -5.49%  '
$ws.Range("D11").Value = "'30.93"
$ws.Range("E11").Value = '  -7.31%  '
$ws.Range("D12").Value = "'0.109"
$ws.Range("E12").Value = '  -2.16%  '
$ws.Range("D13").Value = '2.743.65'
$ws.Range("E13").Value = '  -4.11%  '
$ws.Range("D14").Value = "'6.60"
$ws.Range("E14").Value = '  -4.53%  '
$ws.Range("D15").Value = "'15.06"
$ws.Range("E15").Value = '  -4.99%  '
$ws.Range("D16").Value = '2.347.34'
$ws.Range("E16").Value = '  -5.10%  '
$ws.Range("E17").Value = '  -4.88%  '
$ws.Range("D18").Value = '40.518.08'
$ws.Range("E18").Value = '  -2.95%  '
$ws.Range("D19").Value = '0.0₃0910'
$ws.Range("E19").Value = '  -4.45%  '
$ws.Range("E20").Value = '  -4.96%  '
$ws.Range("D21").Value = "'68.88"
$ws.Range("E21").Value = '  -3.47%  '
$ws.Range("D22").Value = "'10.76"
$ws.Range("E22").Value = '  -5.23%  '
$ws.Range("D23").Value = "'235.70"
$ws.Range("E23").Value = '  -1.86%  '
$ws.Range("D24").Value = "'2.64"
$ws.Range("E24").Value = '  -4.18%  '
$ws.Range("E25").Value = '  +0.45%  '
$ws.Range("E26").Value = '  -6.87%  '
$ws.Range("D27").Value = "'23.53"
$ws.Range("E27").Value = '  -5.25%  '
$ws.Range("E28").Value = '  -3.01%  '
$ws.Range("D29").Value = "'9.37"
$ws.Range("E29").Value = '  -4.83%  '
$ws.Range("D30").Value = "'33.76"
$ws.Range("E30").Value = '  -7.25%  '
$ws.Range("D31").Value = "'155.56"
$ws.Range("E31").Value = '  -1.68%  '
$ws.Range("E32").Value = '  +0.02%  '
$ws.Range("E33").Value = '  -5.86%  '
$ws.Range("D34").Value = "'0.0728"
$ws.Range("E34").Value = '  -5.44%  '
$ws.Range("D35").Value = "'2.42"
$ws.Range("E35").Value = '  -6.67%  '
$ws.Range("E36").Value = '  -2.19%  '
$ws.Range("E37").Value = '  -4.60%  '
$ws.Range("D38").Value = "'15.89"
$ws.Range("E38").Value = '  -8.74%  '
$ws.Range("D39").Value = "'0.0990"
$ws.Range("E39").Value = '  -4.74%  '
$ws.Range("E40").Value = '  -8.78%  '
$ws.Range("E41").Value = '  -5.83%  '
$ws.Range("D42").Value = "'2.27"
$ws.Range("E42").Value = '  -7.91%  '
$ws.Range("D43").Value = '1.961.00'
$ws.Range("E43").Value = '  -1.37%  '
$ws.Range("D44").Value = "'0.0269"
$ws.Range("E44").Value = '  -6.01%  '
$ws.Range("D45").Value = "'17.74"
$ws.Range("E45").Value = '  -8.58%  '
$ws.Range("E46").Value = '  -6.80%  '
$ws.Range("D47").Value = "'9.29"
$ws.Range("E47").Value = '  -0.85%  '
$ws.Range("D48").Value = '2.603.80'
$ws.Range("E48").Value = '  -4.20%  '
$ws.Range("D49").Value = "'93.52"
$ws.Range("E49").Value = '  -4.33%  '
$ws.Range("D50").Value = "'72.72"
$ws.Range("E50").Value = '  -2.43%  '
$ws.Range("D51").Value = "'50.34"
$ws.Range("E51").Value = '  -4.41%  '
